$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows (r2:r9), matching columns A-G; H:M intentionally left blank.
$rows = @(
    @("2024-10-01 13:58:31", "Update", "ABAM",  "RMC", "MRG EP Demo2", "SMH", "Issac Magallanes"),
    @("2024-10-01 13:58:58", "Update", "ABAM",  "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:25:26", "Update", "ABAM",  "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:26:02", "Update", "ABAM",  "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:26:25", "Update", "ABAM2", "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:28:11", "Update", "ABAM2", "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:28:26", "Update", "ABAM2", "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes"),
    @("2024-10-01 14:28:52", "Update", "ABAM2", "RMC", "MRG EP Demo",  "SMH", "Issac Magallanes")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
